# Add data for 2022-05-27
# Updates the "through" date from 05-18 to 05-19, and bumps the May (I6)
# and Total (I14) "2022 (through ...)" column values accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet tab to reflect the new "through" date
$ws.Name = "Through 2022-05-19"

# Update the header label in column I (the shared string cell)
$ws.Range("I1").Value = "2022 (through 05-19)"

# Update the May total (row 6) and the grand Total row (row 14)
$ws.Range("I6").Value = 65
$ws.Range("I14").Value = 617
